$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Permits Filed for 35-17 42nd Street in Astoria, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2025/10/permits-filed-for-35-17-42nd-street-in-astoria-queens.html"
$ws.Range("C2").Value = 'Permits have been filed for a 16-story mixed-use building at 35-17 42nd Street in <a href="https://newyorkyimby.com/neighborhoods/astoria">Astoria</a>, Queens. Located between 75th Street and 76th Street, the lot is one block north of the Freeman Street subway station, served by the 2 and 5 trains. Joel Weiss of Heartfelt Townhouse Build is listed as the owner behind the applications.'
$ws.Range("D2").Value = "2025-10-29T11:00:45+00:00"
$ws.Range("E2").Value = "Wed, 29 Oct 2025 11:00:45 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Astoria"
$ws.Range("H2").Value = ""
